# Adds iMethodInterceptor-style run-manager sheet + a totals sheet, and
# updates selections/active tab to match the authored workbook state.

$wb = $excel.ActiveWorkbook

# --- existing sheet reference ---
$wsTesting = $wb.Worksheets.Item(1)

# --- add "RUNMANAGER" sheet right after "testing" ---
$wsRun = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTesting)
$wsRun.Name = "RUNMANAGER"

# --- add "Sheet1" right after "RUNMANAGER" ---
$wsData = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsRun)
$wsData.Name = "Sheet1"

# ---------------------------------------------------------------------
# RUNMANAGER contents — test run configuration table.
# Cell entry order below mirrors the authored workbook's shared-string
# allocation order (A1, C1, D1, E1, A2, B2, A3, B3, C2, C3, D2, E2, D3,
# E3, then B1 last) so new shared strings land at the same indices.
# ---------------------------------------------------------------------
$wsRun.Range("D1:E1").NumberFormat = "@"
$wsRun.Range("D2:E3").NumberFormat = "@"

$wsRun.Range("A1").Value = "testName"
$wsRun.Range("C1").Value = "execute"
$wsRun.Range("D1").Value = "priority"
$wsRun.Range("E1").Value = "count"

$wsRun.Range("A2").Value = "loginLogoutTest"
$wsRun.Range("B2").Value = "Test the login and logout of the app"

$wsRun.Range("A3").Value = "newTest"
$wsRun.Range("B3").Value = "This is second test"

$wsRun.Range("C2").Value = "no"
$wsRun.Range("C3").Value = "no"

$wsRun.Range("D2").Value = "'1"
$wsRun.Range("E2").Value = "'1"
$wsRun.Range("D3").Value = "'2"
$wsRun.Range("E3").Value = "'1"

$wsRun.Range("B1").Value = "testDescription"

# Column widths roughly matching the authored best-fit sizing.
$wsRun.Columns.Item(1).ColumnWidth = 13.2408854166667
$wsRun.Columns.Item(2).ColumnWidth = 29.7799479166667
$wsRun.Range("D1:E3").Columns.ColumnWidth = 8.39713541666667

$wsRun.PageSetup.PaperSize = 9
$wsRun.PageSetup.Orientation = 1

$wsRun.Activate()
$wsRun.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet1 contents — numeric totals column
# ---------------------------------------------------------------------
$wsData.Range("G6").Value = 125
$wsData.Range("G7").Value = 125
$wsData.Range("G8").Value = 33
$wsData.Range("G9").Value = 41.25
$wsData.Range("G10").Value = 280
$wsData.Range("G11").Value = 76
$wsData.Range("G12").Value = 145
$wsData.Range("G14").Value = 825.25

$wsData.Activate()
$wsData.Range("I16").Select() | Out-Null

# ---------------------------------------------------------------------
# Restore selection on the original sheet, then land back on Sheet1
# (last-activated sheet becomes the active tab, matching activeTab="2").
# ---------------------------------------------------------------------
$wsTesting.Activate()
$wsTesting.Range("H11").Select() | Out-Null

$wsData.Activate()

Write-Output "ok"
